$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.684.33"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.108.54"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'523.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").Value = "'140.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.107.26"
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").Value = "'0.434"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "'7.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").Value = "'0.386"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("D13").Value = "3.640.04"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").Value = "'26.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "57.752.55"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "3.110.46"
$ws.Range("E18").Value = "  +1.56%  "
$ws.Range("D19").Value = "'6.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").Value = "'12.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "'8.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").Value = "'336.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'0.512"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("D25").Value = "'66.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "0.0₃0924"
$ws.Range("E28").Value = "  +3.18%  "
$ws.Range("E29").Value = "  +3.79%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "'7.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'20.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("D35").Value = "'155.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D36").Value = "'4.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.09%  "
$ws.Range("D37").Value = "'6.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("D38").Value = "'27.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("D39").Value = "'1.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("E41").Value = "  +12.92%  "
$ws.Range("D42").Value = "3.150.10"
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("D43").Value = "'0.688"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.08%  "
$ws.Range("D44").Value = "'3.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").Value = "'36.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "2.305.73"
$ws.Range("E47").Value = "  +2.24%  "
$ws.Range("D48").Value = "'0.0259"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("D49").Value = "'0.976"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.22%  "
$ws.Range("D50").Value = "'20.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "'6.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.60%  "
